# Apply the CI_CD说明.docx edit:
#  - paragraphs 2 and 4 lose their <w:pPr> (which only held a
#    <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>) entirely
#  - paragraphs 3 and 5 keep their <w:pPr> (pStyle/numPr/ind) but lose the
#    trailing <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr> inside it
#  - paragraph 6 gains a <w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>
#    (it previously had no <w:pPr> at all)
#  - five new paragraphs are appended at the end of the body (before sectPr)

$d = $word.ActiveDocument

function Wrap-BodyXml([string]$bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Set-ParagraphXml($doc, $index, [string]$innerXml) {
    # InsertXML no-ops when the replacement's visible text is identical to
    # the existing text (it only needs to change paragraph/run properties
    # here). Force a throwaway text change first so the real replacement
    # below is not skipped, then apply the real target XML.
    $paragraph = $doc.Paragraphs.Item($index)
    $range = $paragraph.Range
    $placeholder = Wrap-BodyXml '<w:p><w:r><w:t>@@PLACEHOLDER@@</w:t></w:r></w:p>'
    $range.InsertXML($placeholder)

    $paragraph2 = $doc.Paragraphs.Item($index)
    $range2 = $paragraph2.Range
    $xml = Wrap-BodyXml $innerXml
    $range2.InsertXML($xml)
}

# --- Append five new paragraphs at the very end of the document first ---
# (done before touching paragraph 6's formatting below: InsertXML cannot
# fully swap out the body's very last paragraph mark, so paragraph 6 must
# no longer be the last paragraph before we rewrite it.)
$newParas = '<w:p>' +
        '<w:pPr>' +
            '<w:pStyle w:val="a3"/>' +
            '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
            '<w:ind w:firstLineChars="0"/>' +
            '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>' +
        '</w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>找到反汇编文件</w:t></w:r>' +
    '</w:p>' +
    '<w:p>' +
        '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' +
        '<w:r><w:t>/home/software1/cccc/e902/src</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>/hello</w:t></w:r>' +
        '<w:r><w:t>/</w:t></w:r>' +
        '<w:r><w:t>hello.dis</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>' +
    '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>' +
    '<w:p/>'

$end = $d.Content
$end.Collapse(0)
$xmlEnd = Wrap-BodyXml $newParas
$end.InsertXML($xmlEnd)

# --- Paragraph 2: "cd /home/software1/cccc/e902/src/hello" -------------
# Drop the <w:pPr> that only carried the eastAsia rFonts hint.
$p2 = '<w:p>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>cd</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>/home/software1/cccc/e902/src</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>/hello</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $d 2 $p2

# --- Paragraph 3: "Build hello" -----------------------------------------
# Keep pStyle/numPr/ind, drop the trailing rPr inside pPr.
$p3 = '<w:p>' +
    '<w:pPr>' +
        '<w:pStyle w:val="a3"/>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
        '<w:ind w:firstLineChars="0"/>' +
    '</w:pPr>' +
    '<w:r><w:t>B</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>uild hello</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $d 3 $p3

# --- Paragraph 4: "make all" --------------------------------------------
# Drop the <w:pPr> that only carried the eastAsia rFonts hint.
$p4 = '<w:p>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>make all</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $d 4 $p4

# --- Paragraph 5: "找到bin文件" -------------------------------------------
# Keep pStyle/numPr/ind, drop the trailing rPr inside pPr.
$p5 = '<w:p>' +
    '<w:pPr>' +
        '<w:pStyle w:val="a3"/>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
        '<w:ind w:firstLineChars="0"/>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>找到</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>bin</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>文件</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $d 5 $p5

# --- Paragraph 6: "/home/software1/cccc/e902/bin/hello.bin" -------------
# Gains a <w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>.
$p6 = '<w:p>' +
    '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' +
    '<w:r><w:t>/home/software1/cccc/e902/bin</w:t></w:r>' +
    '<w:r><w:t>/</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>hello.bin</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $d 6 $p6

Write-Output "done"
